$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 1.02
$ws.Range("C2").Value2 = 1.034368082591922
$ws.Range("D2").Value2 = 1.041609298287207
$ws.Range("E2").Value2 = 1.033474052514441
$ws.Range("F2").Value2 = 1.049226066498291
$ws.Range("I2").Value2 = 1.027156265871282
$ws.Range("J2").Value2 = 1.039487710651241
$ws.Range("K2").Value2 = 1.044388328897244
$ws.Range("L2").Value2 = 1.03627629843877
$ws.Range("M2").Value2 = 1.051983715292836
$ws.Range("N2").Value2 = 1.016869072439172
$ws.Range("B3").Value2 = 1.02
$ws.Range("C3").Value2 = 1.036090840985691
$ws.Range("D3").Value2 = 1.043173876150562
$ws.Range("E3").Value2 = 1.034963823625894
$ws.Range("F3").Value2 = 1.050845583589522
$ws.Range("I3").Value2 = 1.027161595642365
$ws.Range("J3").Value2 = 1.040849662557604
$ws.Range("K3").Value2 = 1.045761591950705
$ws.Range("L3").Value2 = 1.03757324102424
$ws.Range("M3").Value2 = 1.053413351595494
$ws.Range("N3").Value2 = 1.017342462170118
$ws.Range("B4").Value2 = 1.02
$ws.Range("C4").Value2 = 1.037203351613426
$ws.Range("D4").Value2 = 1.044184178066086
$ws.Range("E4").Value2 = 1.035926005792827
$ws.Range("F4").Value2 = 1.051890787206321
$ws.Range("I4").Value2 = 1.027162428161306
$ws.Range("J4").Value2 = 1.041728451986604
$ws.Range("K4").Value2 = 1.04664762049356
$ws.Range("L4").Value2 = 1.038410157628912
$ws.Range("M4").Value2 = 1.054335218986712
$ws.Range("N4").Value2 = 1.017647365099939
$ws.Range("B5").Value2 = 1.02
$ws.Range("C5").Value2 = 1.037670530467942
$ws.Range("D5").Value2 = 1.044608420693332
$ws.Range("E5").Value2 = 1.036330086592563
$ws.Range("F5").Value2 = 1.052329547521596
$ws.Range("I5").Value2 = 1.027162151582084
$ws.Range("J5").Value2 = 1.042097310267903
$ws.Range("K5").Value2 = 1.047019502797079
$ws.Range("L5").Value2 = 1.038761458039109
$ws.Range("M5").Value2 = 1.054722015067219
$ws.Range("N5").Value2 = 1.017775211636366
$ws.Range("B6").Value2 = 1.02
$ws.Range("C6").Value2 = 1.037748941702159
$ws.Range("D6").Value2 = 1.044679624573921
$ws.Range("E6").Value2 = 1.036397909168052
$ws.Range("F6").Value2 = 1.052403179899015
$ws.Range("I6").Value2 = 1.027162068393254
$ws.Range("J6").Value2 = 1.04215920919283
$ws.Range("K6").Value2 = 1.047081908332144
$ws.Range("L6").Value2 = 1.038820411567966
$ws.Range("M6").Value2 = 1.054786915725169
$ws.Range("N6").Value2 = 1.017796658104804
$ws.Range("B7").Value2 = 1.02
$ws.Range("C7").Value2 = 1.037209596107824
$ws.Range("D7").Value2 = 1.044189848720653
$ws.Range("E7").Value2 = 1.035931406778784
$ws.Range("F7").Value2 = 1.051896652455808
$ws.Range("I7").Value2 = 1.027162426927843
$ws.Range("J7").Value2 = 1.041733382973045
$ws.Range("K7").Value2 = 1.046652591963177
$ws.Range("L7").Value2 = 1.038414853828443
$ws.Range("M7").Value2 = 1.054340390335638
$ws.Range("N7").Value2 = 1.017649074701923
$ws.Range("B8").Value2 = 1.02
$ws.Range("C8").Value2 = 1.034950766337599
$ws.Range("D8").Value2 = 1.042138492828384
$ws.Range("E8").Value2 = 1.033977906362774
$ws.Range("F8").Value2 = 1.049773962149064
$ws.Range("I8").Value2 = 1.027158608651842
$ws.Range("J8").Value2 = 1.03994850915962
$ws.Range("K8").Value2 = 1.044852966519826
$ws.Range("L8").Value2 = 1.036715086490537
$ws.Range("M8").Value2 = 1.052467536686174
$ws.Range("N8").Value2 = 1.017029351081565
$ws.Range("B9").Value2 = 1.02
$ws.Range("C9").Value2 = 1.030952757457459
$ws.Range("D9").Value2 = 1.038507291430292
$ws.Range("E9").Value2 = 1.030521348962292
$ws.Range("F9").Value2 = 1.046012096935203
$ws.Range("I9").Value2 = 1.027131853018136
$ws.Range("J9").Value2 = 1.036783853914034
$ws.Range("K9").Value2 = 1.041661711735931
$ws.Range("L9").Value2 = 1.033701902996088
$ws.Range("M9").Value2 = 1.049142356187789
$ws.Range("N9").Value2 = 1.01592635418727
$ws.Range("B10").Value2 = 1.02
$ws.Range("C10").Value2 = 1.028274716318853
$ws.Range("D10").Value2 = 1.036074748538535
$ws.Range("E10").Value2 = 1.028206757659134
$ws.Range("F10").Value2 = 1.043489102257777
$ws.Range("I10").Value2 = 1.027100561119203
$ws.Range("J10").Value2 = 1.034660352587585
$ws.Range("K10").Value2 = 1.039520080066466
$ws.Range("L10").Value2 = 1.031680433621443
$ws.Range("M10").Value2 = 1.046908154829334
$ws.Range("N10").Value2 = 1.015183437417094
$ws.Range("B11").Value2 = 1.02
$ws.Range("C11").Value2 = 1.027111897267574
$ws.Range("D11").Value2 = 1.035018492761381
$ws.Range("E11").Value2 = 1.027201945462433
$ws.Range("F11").Value2 = 1.042392885553794
$ws.Range("I11").Value2 = 1.027083821480474
$ws.Range("J11").Value2 = 1.033737454076959
$ws.Range("K11").Value2 = 1.03858923874526
$ws.Range("L11").Value2 = 1.030801974268851
$ws.Range("M11").Value2 = 1.045936450293862
$ws.Range("N11").Value2 = 1.014859899042394
$ws.Range("B12").Value2 = 1.02
$ws.Range("C12").Value2 = 1.026679476273861
$ws.Range("D12").Value2 = 1.034625696064873
$ws.Range("E12").Value2 = 1.026828313471436
$ws.Range("F12").Value2 = 1.041985126033558
$ws.Range("I12").Value2 = 1.027077124378558
$ws.Range("J12").Value2 = 1.033394124054273
$ws.Range("K12").Value2 = 1.038242944718922
$ws.Range("L12").Value2 = 1.030475190564482
$ws.Range("M12").Value2 = 1.045574860627635
$ws.Range("N12").Value2 = 1.014739440271352
$ws.Range("B13").Value2 = 1.02
$ws.Range("C13").Value2 = 1.026772254849608
$ws.Range("D13").Value2 = 1.034709973151465
$ws.Range("E13").Value2 = 1.026908477105743
$ws.Range("F13").Value2 = 1.042072618109457
$ws.Range("I13").Value2 = 1.027078582616981
$ws.Range("J13").Value2 = 1.03346779344054
$ws.Range("K13").Value2 = 1.038317250510816
$ws.Range("L13").Value2 = 1.030545308903105
$ws.Range("M13").Value2 = 1.045652452691611
$ws.Range("N13").Value2 = 1.014765291925261
$ws.Range("B14").Value2 = 1.02
$ws.Range("C14").Value2 = 1.027076163470251
$ws.Range("D14").Value2 = 1.034986033448994
$ws.Range("E14").Value2 = 1.027171069158528
$ws.Range("F14").Value2 = 1.042359191849602
$ws.Range("I14").Value2 = 1.027083277672605
$ws.Range("J14").Value2 = 1.033709085080182
$ws.Range("K14").Value2 = 1.03856062502131
$ws.Range("L14").Value2 = 1.03077497218281
$ws.Range("M14").Value2 = 1.045906574619158
$ws.Range("N14").Value2 = 1.01484994766431
$ws.Range("B15").Value2 = 1.02
$ws.Range("C15").Value2 = 1.027263345099956
$ws.Range("D15").Value2 = 1.035156062449606
$ws.Range("E15").Value2 = 1.027332807406856
$ws.Range("F15").Value2 = 1.042535682681884
$ws.Range("I15").Value2 = 1.027086106945396
$ws.Range("J15").Value2 = 1.033857682960045
$ws.Range("K15").Value2 = 1.038710504417645
$ws.Range("L15").Value2 = 1.030916410721935
$ws.Range("M15").Value2 = 1.046063060319152
$ws.Range("N15").Value2 = 1.0149020693368
$ws.Range("B16").Value2 = 1.02
$ws.Range("C16").Value2 = 1.028351819804038
$ws.Range("D16").Value2 = 1.036144785501769
$ws.Range("E16").Value2 = 1.028273388332793
$ws.Range("F16").Value2 = 1.04356177453904
$ws.Range("I16").Value2 = 1.027101604865038
$ws.Range("J16").Value2 = 1.034721529459916
$ws.Range("K16").Value2 = 1.039581782147816
$ws.Range("L16").Value2 = 1.031738666728805
$ws.Range("M16").Value2 = 1.046972552431582
$ws.Range("N16").Value2 = 1.015204870248428
$ws.Range("B17").Value2 = 1.02
$ws.Range("C17").Value2 = 1.029033720678586
$ws.Range("D17").Value2 = 1.03676418659045
$ws.Range("E17").Value2 = 1.028862690636722
$ws.Range("F17").Value2 = 1.044204403215414
$ws.Range("I17").Value2 = 1.027110472252527
$ws.Range("J17").Value2 = 1.035262476341021
$ws.Range("K17").Value2 = 1.040127365864617
$ws.Range("L17").Value2 = 1.032253594786273
$ws.Range("M17").Value2 = 1.047541898180824
$ws.Range("N17").Value2 = 1.015394310880253
$ws.Range("B18").Value2 = 1.02
$ws.Range("C18").Value2 = 1.029431153336271
$ws.Range("D18").Value2 = 1.037125189167885
$ws.Range("E18").Value2 = 1.029206172730214
$ws.Range("F18").Value2 = 1.044578877385572
$ws.Range("I18").Value2 = 1.027115336658006
$ws.Range("J18").Value2 = 1.035577673482396
$ws.Range("K18").Value2 = 1.040445258681405
$ws.Range("L18").Value2 = 1.032553640491032
$ws.Range("M18").Value2 = 1.047873575741199
$ws.Range("N18").Value2 = 1.015504630134376
$ws.Range("B19").Value2 = 1.02
$ws.Range("C19").Value2 = 1.02956661574261
$ws.Range("D19").Value2 = 1.037248234013572
$ws.Range("E19").Value2 = 1.029323249565262
$ws.Range("F19").Value2 = 1.044706502696974
$ws.Range("I19").Value2 = 1.027116943086824
$ws.Range("J19").Value2 = 1.035685092295276
$ws.Range("K19").Value2 = 1.040553595133076
$ws.Range("L19").Value2 = 1.032655897256173
$ws.Range("M19").Value2 = 1.047986599672073
$ws.Range("N19").Value2 = 1.015542216038983
$ws.Range("B20").Value2 = 1.02
$ws.Range("C20").Value2 = 1.028960591202368
$ws.Range("D20").Value2 = 1.036697760145337
$ws.Range("E20").Value2 = 1.02879948981602
$ws.Range("F20").Value2 = 1.044135492601075
$ws.Range("I20").Value2 = 1.027109552702002
$ws.Range("J20").Value2 = 1.035204471878217
$ws.Range("K20").Value2 = 1.040068864833855
$ws.Range("L20").Value2 = 1.032198379315561
$ws.Range("M20").Value2 = 1.047480855492456
$ws.Range("N20").Value2 = 1.015374004162019
$ws.Range("B21").Value2 = 1.02
$ws.Range("C21").Value2 = 1.026986683839123
$ws.Range("D21").Value2 = 1.034904753274598
$ws.Range("E21").Value2 = 1.027093753487418
$ws.Range("F21").Value2 = 1.042274819030061
$ws.Range("I21").Value2 = 1.027081908325536
$ws.Range("J21").Value2 = 1.033638045282833
$ws.Range("K21").Value2 = 1.038488972246962
$ws.Range("L21").Value2 = 1.030707355547192
$ws.Range("M21").Value2 = 1.045831760234696
$ws.Range("N21").Value2 = 1.014825026483548
$ws.Range("B22").Value2 = 1.02
$ws.Range("C22").Value2 = 1.025742719307434
$ws.Range("D22").Value2 = 1.033774772167689
$ws.Range("E22").Value2 = 1.026018968752129
$ws.Range("F22").Value2 = 1.041101601381904
$ws.Range("I22").Value2 = 1.027061754377649
$ws.Range("J22").Value2 = 1.032650130342138
$ws.Range("K22").Value2 = 1.037492511579905
$ws.Range("L22").Value2 = 1.029767079071645
$ws.Range("M22").Value2 = 1.044791110858421
$ws.Range("N22").Value2 = 1.014478227515966
$ws.Range("B23").Value2 = 1.02
$ws.Range("C23").Value2 = 1.026402447833621
$ws.Range("D23").Value2 = 1.034374051859931
$ws.Range("E23").Value2 = 1.026588956668576
$ws.Range("F23").Value2 = 1.041723867213126
$ws.Range("I23").Value2 = 1.027072701200847
$ws.Range("J23").Value2 = 1.033174135093025
$ws.Range("K23").Value2 = 1.038021054025327
$ws.Range("L23").Value2 = 1.030265807734766
$ws.Range("M23").Value2 = 1.04534314288204
$ws.Range("N23").Value2 = 1.0146622285944
$ws.Range("B24").Value2 = 1.02
$ws.Range("C24").Value2 = 1.028993636218226
$ws.Range("D24").Value2 = 1.03672777627287
$ws.Range("E24").Value2 = 1.028828048312639
$ws.Range("F24").Value2 = 1.044166631454694
$ws.Range("I24").Value2 = 1.027109969158478
$ws.Range("J24").Value2 = 1.035230682611704
$ws.Range("K24").Value2 = 1.04009529997259
$ws.Range("L24").Value2 = 1.032223329746637
$ws.Range("M24").Value2 = 1.047508439329427
$ws.Range("N24").Value2 = 1.015383180444727
$ws.Range("B25").Value2 = 1.02
$ws.Range("C25").Value2 = 1.031988516264757
$ws.Range("D25").Value2 = 1.039448063216481
$ws.Range("E25").Value2 = 1.031416705469276
$ws.Range("F25").Value2 = 1.046987237175486
$ws.Range("I25").Value2 = 1.027141142171841
$ws.Range("J25").Value2 = 1.037604366462381
$ws.Range("K25").Value2 = 1.042489171637353
$ws.Range("L25").Value2 = 1.034483074221182
$ws.Range("M25").Value2 = 1.050005018279794
$ws.Range("N25").Value2 = 1.016212826519011
